$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the values from A1:A4 into A5:A8 (same sequence repeated)
$ws.Range("A5").Value = 3174466432
$ws.Range("A6").Value = 3247439861
$ws.Range("A7").Value = 3104023154
$ws.Range("A8").Value = 3215996243

# Update the selection to match the new active range
$ws.Range("A5:A8").Select()
